$d = $word.ActiveDocument

# 1. Update the letter date: September 19, 2025 -> September 21, 2025
$null = $d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, $true, 1, $false, "September 21, 2025", 2)

# 2. Split the sender's mailing address paragraph "969 Story Road, San Jose CA 95122"
#    into two paragraphs: "969 Story Road" and "San Jose, CA 95122".
#    (Only the sender address near the top of the letter changes - NOT the
#    "PROPERTY ADDRESS:" table cell further down, which keeps the original text.)
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*969 Story Road, San Jose CA 95122*" -and $p.Range.Information(12) -eq $false) {
        $targetPara = $p
        break
    }
}
if ($targetPara -ne $null) {
    $full = $targetPara.Range
    $origEnd = $full.End
    $splitOffset = $full.Start + ("969 Story Road").Length
    # Collapse a zero-length range right after "969 Story Road" and insert
    # a paragraph break there - this duplicates the run formatting
    # (Arial, sz/szCs 22) onto both sides of the split.
    $collapse = $d.Range($splitOffset, $splitOffset)
    $collapse.InsertParagraphAfter()
    # The inserted paragraph mark occupies $splitOffset, so everything that
    # used to start there (and the original paragraph end) is now shifted by
    # one character.
    $tail = $d.Range($splitOffset + 1, $origEnd + 1)
    $tail.Text = "San Jose, CA 95122"
}

# 3. Remove the now-redundant empty "No Spacing" paragraph that sat right
#    after the "Board of Directors" line.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -like "*Board of Directors*") {
        $next = $p.Next()
        if ($next.Range.Text -eq "`r" -and $next.Style.NameLocal -eq "No Spacing") {
            $next.Range.Delete()
        }
        break
    }
}
